$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 961.4706
$ws.Range("I19").Value = 806.55554
$ws.Range("J19").Value = 1135.75
$ws.Range("K19").Value = 806.55554
$ws.Range("L19").Value = 1135.75
$ws.Range("M19").Value = -631.55554
$ws.Range("N19").Value = -1485.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 462.78125
$ws.Range("I33").Value = 160.24138
$ws.Range("J33").Value = 3387.3333
$ws.Range("K33").Value = 160.24138
$ws.Range("L33").Value = 3387.3333
$ws.Range("M33").Value = 68.75862000000001
$ws.Range("N33").Value = -3845.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 889.08
$ws.Range("I129").Value = 595.5
$ws.Range("J129").Value = 914.6087
$ws.Range("K129").Value = 1786.5
$ws.Range("L129").Value = 2743.8261
$ws.Range("M129").Value = 3213.5
$ws.Range("N129").Value = -12743.8261

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 5106978
$ws.Range("I132").Value = 5560778.5
$ws.Range("J132").Value = 1724.75
$ws.Range("K132").Value = 16682335.5
$ws.Range("L132").Value = 5174.25
$ws.Range("M132").Value = -16679805.5
$ws.Range("N132").Value = -10234.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 2505.5715
$ws.Range("I135").Value = 669.2857
$ws.Range("J135").Value = 4341.857
$ws.Range("K135").Value = 6023.571300000001
$ws.Range("L135").Value = 39076.713
$ws.Range("M135").Value = -3488.571300000001
$ws.Range("N135").Value = -44146.713

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1576.4
$ws.Range("I137").Value = 1457.1305
$ws.Range("J137").Value = 1805
$ws.Range("K137").Value = 4371.3915
$ws.Range("L137").Value = 5415
$ws.Range("M137").Value = -1821.3915
$ws.Range("N137").Value = -10515

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2821.6
$ws.Range("I138").Value = 3229.25
$ws.Range("J138").Value = 2549.8333
$ws.Range("K138").Value = 9687.75
$ws.Range("L138").Value = 7649.499899999999
$ws.Range("M138").Value = -4547.75
$ws.Range("N138").Value = -17929.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 2000
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2000
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -814
$ws.Range("N86").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1259.6
$ws.Range("I88").Value = 1274.5
$ws.Range("J88").Value = 1200
$ws.Range("K88").Value = 1274.5
$ws.Range("L88").Value = 1200
$ws.Range("M88").Value = -868.5
$ws.Range("N88").Value = -2012

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 2000
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 6000
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -72
$ws.Range("N89").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 1259.6
$ws.Range("I91").Value = 1274.5
$ws.Range("J91").Value = 1200
$ws.Range("K91").Value = 1274.5
$ws.Range("L91").Value = 1200
$ws.Range("M91").Value = 129.5
$ws.Range("N91").Value = -4008

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 24950
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 24950
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 24950
$ws.Range("N11").Value = -25230
$ws.Range("M11").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 67679.3
$ws.Range("I86").Value = 126310.22
$ws.Range("J86").Value = 1719.5
$ws.Range("K86").Value = 126310.22
$ws.Range("L86").Value = 1719.5
$ws.Range("M86").Value = -125187.22
$ws.Range("N86").Value = -3965.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 67679.3
$ws.Range("I89").Value = 126310.22
$ws.Range("J89").Value = 1719.5
$ws.Range("K89").Value = 631551.1
$ws.Range("L89").Value = 8597.5
$ws.Range("M89").Value = -625935.1
$ws.Range("N89").Value = -19829.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1614.6296
$ws.Range("I99").Value = 1364.2
$ws.Range("J99").Value = 1671.5454
$ws.Range("K99").Value = 1364.2
$ws.Range("L99").Value = 1671.5454
$ws.Range("M99").Value = 133.8
$ws.Range("N99").Value = -4667.5454

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2252.7302
$ws.Range("I31").Value = 1436.7931
$ws.Range("J31").Value = 2948.6765
$ws.Range("K31").Value = 1436.7931
$ws.Range("L31").Value = 2948.6765
$ws.Range("M31").Value = -1141.7931
$ws.Range("N31").Value = -3538.6765

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2252.7302
$ws.Range("I34").Value = 1436.7931
$ws.Range("J34").Value = 2948.6765
$ws.Range("K34").Value = 1436.7931
$ws.Range("L34").Value = 2948.6765
$ws.Range("M34").Value = -1234.7931
$ws.Range("N34").Value = -3352.6765

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5343.857
$ws.Range("I132").Value = 5916.846
$ws.Range("J132").Value = 4412.75
$ws.Range("K132").Value = 17750.538
$ws.Range("L132").Value = 13238.25
$ws.Range("M132").Value = -15220.538

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1429.5714
$ws.Range("I5").Value = 1514.2
$ws.Range("J5").Value = 1352.6364
$ws.Range("K5").Value = 4542.6
$ws.Range("L5").Value = 4057.9092
$ws.Range("M5").Value = -4430.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1288.7142
$ws.Range("I113").Value = 2474.4
$ws.Range("J113").Value = 630
$ws.Range("K113").Value = 7423.200000000001
$ws.Range("L113").Value = 1890
$ws.Range("M113").Value = -5253.200000000001
$ws.Range("N113").Value = -6230

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1540.375
$ws.Range("I132").Value = 695.2632
$ws.Range("J132").Value = 2305
$ws.Range("K132").Value = 6257.3688
$ws.Range("L132").Value = 20745
$ws.Range("M132").Value = -3727.3688
$ws.Range("N132").Value = -25805

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1429.5714
$ws.Range("I135").Value = 1514.2
$ws.Range("J135").Value = 1352.6364
$ws.Range("K135").Value = 13627.8
$ws.Range("L135").Value = 12173.7276
$ws.Range("M135").Value = -11092.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 2001.2
$ws.Range("I139").Value = 864.53845
$ws.Range("J139").Value = 3232.5833
$ws.Range("K139").Value = 2593.61535
$ws.Range("L139").Value = 9697.749899999999
$ws.Range("M139").Value = 2546.38465

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 100105016
$ws.Range("I80").Value = 200207400
$ws.Range("J80").Value = 2635.8
$ws.Range("K80").Value = 200207400
$ws.Range("L80").Value = 2635.8
$ws.Range("M80").Value = -200206402
$ws.Range("N80").Value = -4631.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 100105016
$ws.Range("I83").Value = 200207400
$ws.Range("J83").Value = 2635.8
$ws.Range("K83").Value = 1001037000
$ws.Range("L83").Value = 13179
$ws.Range("M83").Value = -1001032008
$ws.Range("N83").Value = -23163

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 15000
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 15000
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 15000
$ws.Range("N96").Value = -20492

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1499.909
$ws.Range("I113").Value = 900
$ws.Range("J113").Value = 1633.2222
$ws.Range("K113").Value = 900
$ws.Range("L113").Value = 1633.2222
$ws.Range("M113").Value = 1270
$ws.Range("N113").Value = -5973.2222

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3645.7693
$ws.Range("I68").Value = 2277
$ws.Range("J68").Value = 4254.1113
$ws.Range("K68").Value = 2277
$ws.Range("L68").Value = 4254.1113
$ws.Range("M68").Value = -1528
$ws.Range("N68").Value = -5752.1113

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3645.7693
$ws.Range("I71").Value = 2277
$ws.Range("J71").Value = 4254.1113
$ws.Range("K71").Value = 11385
$ws.Range("L71").Value = 21270.5565
$ws.Range("M71").Value = -7641
$ws.Range("N71").Value = -28758.5565

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1827.8
$ws.Range("I100").Value = 1199.8334
$ws.Range("J100").Value = 2769.75
$ws.Range("K100").Value = 1199.8334
$ws.Range("L100").Value = 2769.75
$ws.Range("M100").Value = -658.8334
$ws.Range("N100").Value = -3851.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H135").Value = 39080
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 39080
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 39080
$ws.Range("N135").Value = -49220

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H137").Value = 35333.332
$ws.Range("I137").Value = 1000
$ws.Range("J137").Value = 42200
$ws.Range("K137").Value = 1000
$ws.Range("L137").Value = 42200
$ws.Range("M137").Value = 4100
$ws.Range("N137").Value = -52400

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 65712.5
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 65712.5
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 65712.5
$ws.Range("N139").Value = -75992.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2294.2632
$ws.Range("I122").Value = 1613.8572
$ws.Range("J122").Value = 2691.1667
$ws.Range("K122").Value = 4841.571599999999
$ws.Range("L122").Value = 8073.500100000001
$ws.Range("M122").Value = -2391.571599999999
$ws.Range("N122").Value = -12973.5001
